# Auto-generated Excel COM-interop edit script
# Applies per-cell numeric updates to sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(74, 8).Value = 7887  # H74: was 6943.25
$ws.Cells.Item(74, 9).Value = 5770  # I74: was 5923
$ws.Cells.Item(74, 11).Value = 5770  # K74: was 5923
$ws.Cells.Item(74, 13).Value = -4834  # M74: was -4987

$ws.Cells.Item(76, 8).Value = 5070.8  # H76: was 4817.1665
$ws.Cells.Item(76, 9).Value = 5450  # I76: was 4816.3335
$ws.Cells.Item(76, 11).Value = 5450  # K76: was 4816.3335
$ws.Cells.Item(76, 13).Value = -5135  # M76: was -4501.3335

$ws.Cells.Item(77, 8).Value = 7887  # H77: was 6943.25
$ws.Cells.Item(77, 9).Value = 5770  # I77: was 5923
$ws.Cells.Item(77, 11).Value = 28850  # K77: was 29615
$ws.Cells.Item(77, 13).Value = -24170  # M77: was -24935

$ws.Cells.Item(79, 8).Value = 5070.8  # H79: was 4817.1665
$ws.Cells.Item(79, 9).Value = 5450  # I79: was 4816.3335
$ws.Cells.Item(79, 11).Value = 5450  # K79: was 4816.3335
$ws.Cells.Item(79, 13).Value = -4358  # M79: was -3724.3335

$ws.Cells.Item(107, 8).Value = 750.9091  # H107: was 749.26086
$ws.Cells.Item(107, 9).Value = 750.9091  # I107: was 749.26086
$ws.Cells.Item(107, 11).Value = 750.9091  # K107: was 749.26086
$ws.Cells.Item(107, 13).Value = 1169.0909  # M107: was 1170.73914

$ws.Cells.Item(116, 8).Value = 69035.336  # H116: was 73731
$ws.Cells.Item(116, 9).Value = 127066.25  # I116: was 144747.72
$ws.Cells.Item(116, 11).Value = 127066.25  # K116: was 144747.72
$ws.Cells.Item(116, 13).Value = -123624.25  # M116: was -141305.72

$ws.Cells.Item(132, 8).Value = 53858.79  # H132: was 56869.89
$ws.Cells.Item(132, 9).Value = 56738.332  # I132: was 60095.94
$ws.Cells.Item(132, 11).Value = 170214.996  # K132: was 180287.82
$ws.Cells.Item(132, 13).Value = -167684.996  # M132: was -177757.82

$ws.Cells.Item(135, 8).Value = 626.9474  # H135: was 649.775
$ws.Cells.Item(135, 9).Value = 450.65625  # I135: was 439.24243
$ws.Cells.Item(135, 10).Value = 1567.1666  # J135: was 1642.2858
$ws.Cells.Item(135, 11).Value = 4055.90625  # K135: was 3953.18187
$ws.Cells.Item(135, 12).Value = 14104.4994  # L135: was 14780.5722
$ws.Cells.Item(135, 13).Value = -1520.90625  # M135: was -1418.18187
$ws.Cells.Item(135, 14).Value = -19174.4994  # N135: was -19850.5722

$ws.Cells.Item(137, 8).Value = 4514.75  # H137: was 7783.6665
$ws.Cells.Item(137, 9).Value = 2521.75  # I137: was 4899.6665
$ws.Cells.Item(137, 10).Value = 8500.75  # J137: was 10667.667
$ws.Cells.Item(137, 11).Value = 7565.25  # K137: was 14698.9995
$ws.Cells.Item(137, 12).Value = 25502.25  # L137: was 32003.001
$ws.Cells.Item(137, 13).Value = -5015.25  # M137: was -12148.9995
$ws.Cells.Item(137, 14).Value = -30602.25  # N137: was -37103.001

$ws.Cells.Item(138, 8).Value = 2368.8667  # H138: was 2593.4119
$ws.Cells.Item(138, 9).Value = 2012.2727  # I138: was 2064
$ws.Cells.Item(138, 10).Value = 3349.5  # J138: was 3349.7144
$ws.Cells.Item(138, 11).Value = 6036.8181  # K138: was 6192
$ws.Cells.Item(138, 12).Value = 10048.5  # L138: was 10049.1432
$ws.Cells.Item(138, 13).Value = -896.8181000000004  # M138: was -1052
$ws.Cells.Item(138, 14).Value = -20328.5  # N138: was -20329.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2958.9546  # H32: was 3200.158
$ws.Cells.Item(32, 9).Value = 3004.6191  # I32: was 3200.158
$ws.Cells.Item(32, 10).Value = 2000  # J32: was 0
$ws.Cells.Item(32, 11).Value = 3004.6191  # K32: was 3200.158
$ws.Cells.Item(32, 12).Value = 2000  # L32: was 0
$ws.Cells.Item(32, 13).Value = -2717.6191  # M32: was -2913.158
$ws.Cells.Item(32, 14).Value = -2574  # N32: was None

$ws.Cells.Item(74, 8).Value = 3859239.5  # H74: was 4210017
$ws.Cells.Item(74, 9).Value = 1853240.8  # I74: was 2059080
$ws.Cells.Item(74, 11).Value = 1853240.8  # K74: was 2059080
$ws.Cells.Item(74, 13).Value = -1852366.8  # M74: was -2058206

$ws.Cells.Item(77, 8).Value = 3859239.5  # H77: was 4210017
$ws.Cells.Item(77, 9).Value = 1853240.8  # I77: was 2059080
$ws.Cells.Item(77, 11).Value = 9266204  # K77: was 10295400
$ws.Cells.Item(77, 13).Value = -9261836  # M77: was -10291032

$ws.Cells.Item(88, 8).Value = 3000  # H88: was 1733.4667
$ws.Cells.Item(88, 9).Value = 3000  # I88: was 1869.6666
$ws.Cells.Item(88, 10).Value = 3000  # J88: was 1529.1666
$ws.Cells.Item(88, 11).Value = 3000  # K88: was 1869.6666
$ws.Cells.Item(88, 12).Value = 3000  # L88: was 1529.1666
$ws.Cells.Item(88, 13).Value = -2594  # M88: was -1463.6666
$ws.Cells.Item(88, 14).Value = -3812  # N88: was -2341.1666

$ws.Cells.Item(91, 8).Value = 3000  # H91: was 1733.4667
$ws.Cells.Item(91, 9).Value = 3000  # I91: was 1869.6666
$ws.Cells.Item(91, 10).Value = 3000  # J91: was 1529.1666
$ws.Cells.Item(91, 11).Value = 3000  # K91: was 1869.6666
$ws.Cells.Item(91, 12).Value = 3000  # L91: was 1529.1666
$ws.Cells.Item(91, 13).Value = -1596  # M91: was -465.6666
$ws.Cells.Item(91, 14).Value = -5808  # N91: was -4337.1666

$ws.Cells.Item(132, 8).Value = 23811954  # H132: was 26318272
$ws.Cells.Item(132, 9).Value = 2503  # I132: was 2571.9375
$ws.Cells.Item(132, 11).Value = 7509  # K132: was 7715.8125
$ws.Cells.Item(132, 13).Value = -4979  # M132: was -5185.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(19, 8).Value = 0  # H19: was 487
$ws.Cells.Item(19, 9).Value = 0  # I19: was 487
$ws.Cells.Item(19, 11).Value = 0  # K19: was 487
$ws.Cells.Item(19, 13).ClearContents()  # M19: was -314

$ws.Cells.Item(20, 8).Value = 1000  # H20: was 0
$ws.Cells.Item(20, 9).Value = 1000  # I20: was 0
$ws.Cells.Item(20, 11).Value = 1000  # K20: was 0
$ws.Cells.Item(20, 13).Value = -753  # M20: was None

$ws.Cells.Item(35, 8).Value = 42000  # H35: was 41995
$ws.Cells.Item(35, 10).Value = 42000  # J35: was 41995
$ws.Cells.Item(35, 12).Value = 42000  # L35: was 41995
$ws.Cells.Item(35, 14).Value = -42620  # N35: was -42615

$ws.Cells.Item(86, 8).Value = 3006.818  # H86: was 2080.0435
$ws.Cells.Item(86, 9).Value = 2318.5  # I86: was 1561.3334
$ws.Cells.Item(86, 10).Value = 3400.1428  # J86: was 2645.9092
$ws.Cells.Item(86, 11).Value = 2318.5  # K86: was 1561.3334
$ws.Cells.Item(86, 12).Value = 3400.1428  # L86: was 2645.9092
$ws.Cells.Item(86, 13).Value = -1195.5  # M86: was -438.3334
$ws.Cells.Item(86, 14).Value = -5646.1428  # N86: was -4891.9092

$ws.Cells.Item(89, 8).Value = 3006.818  # H89: was 2080.0435
$ws.Cells.Item(89, 9).Value = 2318.5  # I89: was 1561.3334
$ws.Cells.Item(89, 10).Value = 3400.1428  # J89: was 2645.9092
$ws.Cells.Item(89, 11).Value = 11592.5  # K89: was 7806.666999999999
$ws.Cells.Item(89, 12).Value = 17000.714  # L89: was 13229.546
$ws.Cells.Item(89, 13).Value = -5976.5  # M89: was -2190.666999999999
$ws.Cells.Item(89, 14).Value = -28232.714  # N89: was -24461.546

$ws.Cells.Item(105, 8).Value = 3565.6667  # H105: was 3598.5
$ws.Cells.Item(105, 9).Value = 3565.6667  # I105: was 3598.5
$ws.Cells.Item(105, 11).Value = 3565.6667  # K105: was 3598.5
$ws.Cells.Item(105, 13).Value = -1818.6667  # M105: was -1851.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 10188.8  # H62: was 8906.833000000001
$ws.Cells.Item(62, 9).Value = 2736  # I62: was 2688.2
$ws.Cells.Item(62, 11).Value = 2736  # K62: was 2688.2
$ws.Cells.Item(62, 13).Value = -2112  # M62: was -2064.2

$ws.Cells.Item(65, 8).Value = 10188.8  # H65: was 8906.833000000001
$ws.Cells.Item(65, 9).Value = 2736  # I65: was 2688.2
$ws.Cells.Item(65, 11).Value = 13680  # K65: was 13441
$ws.Cells.Item(65, 13).Value = -10560  # M65: was -10321

$ws.Cells.Item(132, 8).Value = 7052.08  # H132: was 7490.8096
$ws.Cells.Item(132, 9).Value = 7259.1055  # I132: was 7642.706
$ws.Cells.Item(132, 10).Value = 6396.5  # J132: was 6845.25
$ws.Cells.Item(132, 11).Value = 21777.3165  # K132: was 22928.118
$ws.Cells.Item(132, 12).Value = 19189.5  # L132: was 20535.75
$ws.Cells.Item(132, 13).Value = -19247.3165  # M132: was -20398.118
$ws.Cells.Item(132, 14).Value = -24249.5  # N132: was -25595.75

$ws.Cells.Item(134, 8).Value = 7146779.5  # H134: was 8337960
$ws.Cells.Item(134, 9).Value = 2555.5715  # I134: was 3698.4
$ws.Cells.Item(134, 11).Value = 7666.7145  # K134: was 11095.2
$ws.Cells.Item(134, 13).Value = -5131.7145  # M134: was -8560.200000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(18, 8).Value = 2044.4445  # H18: was 2056.7778
$ws.Cells.Item(18, 9).Value = 1486.2858  # I18: was 1502.1428
$ws.Cells.Item(18, 11).Value = 4458.857400000001  # K18: was 4506.428400000001
$ws.Cells.Item(18, 13).Value = -4289.857400000001  # M18: was -4337.428400000001

$ws.Cells.Item(37, 8).Value = 119998.2  # H37: was 119998.5
$ws.Cells.Item(37, 10).Value = 119998.2  # J37: was 119998.5
$ws.Cells.Item(37, 12).Value = 359994.6  # L37: was 359995.5
$ws.Cells.Item(37, 14).Value = -360218.6  # N37: was -360219.5

$ws.Cells.Item(98, 8).Value = 399.125  # H98: was 402.9
$ws.Cells.Item(98, 10).Value = 800  # J98: was 545.3333
$ws.Cells.Item(98, 12).Value = 2400  # L98: was 1635.9999
$ws.Cells.Item(98, 14).Value = -5396  # N98: was -4631.9999

$ws.Cells.Item(103, 8).Value = 255.27272  # H103: was 255.45454
$ws.Cells.Item(103, 9).Value = 283  # I103: was 266.7143
$ws.Cells.Item(103, 10).Value = 222  # J103: was 235.75
$ws.Cells.Item(103, 11).Value = 849  # K103: was 800.1428999999999
$ws.Cells.Item(103, 12).Value = 666  # L103: was 707.25
$ws.Cells.Item(103, 13).Value = 30  # M103: was 78.85710000000006
$ws.Cells.Item(103, 14).Value = -2424  # N103: was -2465.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 12166.667  # H70: was 10484.923
$ws.Cells.Item(70, 9).Value = 11916.667  # I70: was 9630.5
$ws.Cells.Item(70, 10).Value = 12666.667  # J70: was 13333
$ws.Cells.Item(70, 11).Value = 11916.667  # K70: was 9630.5
$ws.Cells.Item(70, 12).Value = 12666.667  # L70: was 13333
$ws.Cells.Item(70, 13).Value = -11646.667  # M70: was -9360.5
$ws.Cells.Item(70, 14).Value = -13206.667  # N70: was -13873

$ws.Cells.Item(73, 8).Value = 12166.667  # H73: was 10484.923
$ws.Cells.Item(73, 9).Value = 11916.667  # I73: was 9630.5
$ws.Cells.Item(73, 10).Value = 12666.667  # J73: was 13333
$ws.Cells.Item(73, 11).Value = 11916.667  # K73: was 9630.5
$ws.Cells.Item(73, 12).Value = 12666.667  # L73: was 13333
$ws.Cells.Item(73, 13).Value = -10980.667  # M73: was -8694.5
$ws.Cells.Item(73, 14).Value = -14538.667  # N73: was -15205

$ws.Cells.Item(80, 8).Value = 7391.3335  # H80: was 4576.5454
$ws.Cells.Item(80, 9).Value = 3869.6  # I80: was 2867.875
$ws.Cells.Item(80, 10).Value = 25000  # J80: was 9133
$ws.Cells.Item(80, 11).Value = 3869.6  # K80: was 2867.875
$ws.Cells.Item(80, 12).Value = 25000  # L80: was 9133
$ws.Cells.Item(80, 13).Value = -2871.6  # M80: was -1869.875
$ws.Cells.Item(80, 14).Value = -26996  # N80: was -11129

$ws.Cells.Item(83, 8).Value = 7391.3335  # H83: was 4576.5454
$ws.Cells.Item(83, 9).Value = 3869.6  # I83: was 2867.875
$ws.Cells.Item(83, 10).Value = 25000  # J83: was 9133
$ws.Cells.Item(83, 11).Value = 19348  # K83: was 14339.375
$ws.Cells.Item(83, 12).Value = 125000  # L83: was 45665
$ws.Cells.Item(83, 13).Value = -14356  # M83: was -9347.375
$ws.Cells.Item(83, 14).Value = -134984  # N83: was -55649

$ws.Cells.Item(126, 8).Value = 5174.467  # H126: was 4983.875
$ws.Cells.Item(126, 9).Value = 10530.75  # I126: was 8849.6
$ws.Cells.Item(126, 11).Value = 31592.25  # K126: was 26548.8
$ws.Cells.Item(126, 13).Value = -29122.25  # M126: was -24078.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(23, 8).Value = 5953  # H23: was 4501.5
$ws.Cells.Item(23, 9).Value = 5953  # I23: was 4501.5
$ws.Cells.Item(23, 11).Value = 5953  # K23: was 4501.5
$ws.Cells.Item(23, 13).Value = -5723  # M23: was -4271.5

$ws.Cells.Item(36, 8).Value = 0  # H36: was 94996.5
$ws.Cells.Item(36, 10).Value = 0  # J36: was 94996.5
$ws.Cells.Item(36, 12).Value = 0  # L36: was 94996.5
$ws.Cells.Item(36, 14).ClearContents()  # N36: was -96120.5

$ws.Cells.Item(46, 8).Value = 2593.4849  # H46: was 2617.7273
$ws.Cells.Item(46, 9).Value = 1500  # I46: was 1900
$ws.Cells.Item(46, 10).Value = 2744.3103  # J46: was 2689.5
$ws.Cells.Item(46, 11).Value = 1500  # K46: was 1900
$ws.Cells.Item(46, 12).Value = 2744.3103  # L46: was 2689.5
$ws.Cells.Item(46, 13).Value = -1312  # M46: was -1712
$ws.Cells.Item(46, 14).Value = -3120.3103  # N46: was -3065.5

$ws.Cells.Item(122, 8).Value = 3280.5454  # H122: was 3409.6667
$ws.Cells.Item(122, 9).Value = 3000.25  # I122: was 3167.6667
$ws.Cells.Item(122, 10).Value = 3616.9  # J122: was 3651.6667
$ws.Cells.Item(122, 11).Value = 9000.75  # K122: was 9503.000100000001
$ws.Cells.Item(122, 12).Value = 10850.7  # L122: was 10955.0001
$ws.Cells.Item(122, 13).Value = -6550.75  # M122: was -7053.000100000001
$ws.Cells.Item(122, 14).Value = -15750.7  # N122: was -15855.0001

$ws.Cells.Item(132, 8).Value = 0  # H132: was 2190.923
$ws.Cells.Item(132, 9).Value = 0  # I132: was 1956.8334
$ws.Cells.Item(132, 10).Value = 0  # J132: was 5000
$ws.Cells.Item(132, 11).Value = 0  # K132: was 5870.5002
$ws.Cells.Item(132, 12).Value = 0  # L132: was 15000
$ws.Cells.Item(132, 13).ClearContents()  # M132: was -3340.5002
$ws.Cells.Item(132, 14).ClearContents()  # N132: was -20060

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(23, 8).Value = 2597  # H23: was 3112.6
$ws.Cells.Item(23, 9).Value = 2669.6365  # I23: was 3258.6667
$ws.Cells.Item(23, 11).Value = 2669.6365  # K23: was 3258.6667
$ws.Cells.Item(23, 13).Value = -2440.6365  # M23: was -3029.6667

$ws.Cells.Item(29, 8).Value = 1950  # H29: was 1933
$ws.Cells.Item(29, 9).Value = 1950  # I29: was 1933
$ws.Cells.Item(29, 11).Value = 1950  # K29: was 1933
$ws.Cells.Item(29, 13).Value = -1660  # M29: was -1643

$ws.Cells.Item(41, 8).Value = 83993.75  # H41: was 83994
$ws.Cells.Item(41, 10).Value = 91055.10000000001  # J41: was 91055.39999999999
$ws.Cells.Item(41, 12).Value = 91055.10000000001  # L41: was 91055.39999999999
$ws.Cells.Item(41, 14).Value = -91835.10000000001  # N41: was -91835.39999999999

$ws.Cells.Item(44, 8).Value = 14000  # H44: was 23500
$ws.Cells.Item(44, 10).Value = 14000  # J44: was 23500
$ws.Cells.Item(44, 12).Value = 14000  # L44: was 23500
$ws.Cells.Item(44, 14).Value = -15108  # N44: was -24608

$ws.Cells.Item(50, 8).Value = 0  # H50: was 21000
$ws.Cells.Item(50, 9).Value = 0  # I50: was 21000
$ws.Cells.Item(50, 11).Value = 0  # K50: was 21000
$ws.Cells.Item(50, 13).ClearContents()  # M50: was -20369

$ws.Cells.Item(81, 8).Value = 2780.818  # H81: was 2447.0715
$ws.Cells.Item(81, 9).Value = 1835.875  # I81: was 1688.0834
$ws.Cells.Item(81, 10).Value = 5300.6665  # J81: was 7001
$ws.Cells.Item(81, 11).Value = 3671.75  # K81: was 3376.1668
$ws.Cells.Item(81, 12).Value = 10601.333  # L81: was 14002
$ws.Cells.Item(81, 13).Value = -2610.75  # M81: was -2315.1668
$ws.Cells.Item(81, 14).Value = -12723.333  # N81: was -16124

$ws.Cells.Item(84, 8).Value = 2780.818  # H84: was 2447.0715
$ws.Cells.Item(84, 9).Value = 1835.875  # I84: was 1688.0834
$ws.Cells.Item(84, 10).Value = 5300.6665  # J84: was 7001
$ws.Cells.Item(84, 11).Value = 18358.75  # K84: was 16880.834
$ws.Cells.Item(84, 12).Value = 53006.665  # L84: was 70010
$ws.Cells.Item(84, 13).Value = -13054.75  # M84: was -11576.834
$ws.Cells.Item(84, 14).Value = -63614.665  # N84: was -80618

$ws.Cells.Item(107, 8).Value = 1192.8695  # H107: was 1314.9048
$ws.Cells.Item(107, 9).Value = 734.5294  # I107: was 837.4666999999999
$ws.Cells.Item(107, 10).Value = 2491.5  # J107: was 2508.5
$ws.Cells.Item(107, 11).Value = 2203.5882  # K107: was 2512.4001
$ws.Cells.Item(107, 12).Value = 7474.5  # L107: was 7525.5
$ws.Cells.Item(107, 13).Value = -283.5882000000001  # M107: was -592.4000999999998
$ws.Cells.Item(107, 14).Value = -11314.5  # N107: was -11365.5

$ws.Cells.Item(122, 8).Value = 3819.6875  # H122: was 3993.7144
$ws.Cells.Item(122, 9).Value = 3741  # I122: was 3916.3076
$ws.Cells.Item(122, 11).Value = 11223  # K122: was 11748.9228
$ws.Cells.Item(122, 13).Value = -8773  # M122: was -9298.9228

$ws.Cells.Item(132, 8).Value = 1543.6364  # H132: was 1635.75
$ws.Cells.Item(132, 9).Value = 1543.6364  # I132: was 1635.75
$ws.Cells.Item(132, 11).Value = 4630.9092  # K132: was 4907.25
$ws.Cells.Item(132, 13).Value = -2100.9092  # M132: was -2377.25
